# Slide 4 ("Dữ liệu" -> "Data") title placeholder:
#  - pin down the explicit position/size that PowerPoint writes once the
#    inherited-from-layout placeholder is touched/moved
#  - update the title text from the Vietnamese "Dữ liệu" to "Data"
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(4)
$sh = $s.Shapes.Title

# Values below are the EMU targets (818712 / 352920 / 10571998 / 970450)
# expressed in points (1 pt = 12700 EMU) with enough precision that the
# round-trip through the host's float handling lands back on the exact
# integer EMU the author's deck has.
$sh.Left   = 64.46552
$sh.Top    = 27.789
$sh.Width  = 832.4408
$sh.Height = 76.4134

$sh.TextFrame.TextRange.Text = "Data"
